# Adds the new match row (row 77) to the Thai League 1 2023-2024 sheet,
# mirroring the formatting of the last existing data row (76).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles) from the last existing row so the new row keeps
# the same look (bold/bordered/centered index column, date-formatted column E).
$ws.Range("A76:V76").Copy()
$ws.Range("A77:V77").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new match data.
$ws.Range("A77").Value = 76
$ws.Range("B77").Value = "thailand"
$ws.Range("C77").Value = "thai-league-1"
$ws.Range("D77").Value = "2023-2024"
$ws.Range("E77").Value = 45240.5
$ws.Range("F77").Value = "Sukhothai"
$ws.Range("G77").Value = 3
$ws.Range("H77").Value = "Chonburi"
$ws.Range("I77").Value = 2
$ws.Range("J77").Value = 2.81
$ws.Range("K77").Value = "05/11/2023 13:13"
$ws.Range("L77").Value = 2.85
$ws.Range("M77").Value = "10/11/2023 11:59"
$ws.Range("N77").Value = 3.42
$ws.Range("O77").Value = "05/11/2023 13:13"
$ws.Range("P77").Value = 3.58
$ws.Range("Q77").Value = "10/11/2023 11:51"
$ws.Range("R77").Value = 2.47
$ws.Range("S77").Value = "05/11/2023 13:13"
$ws.Range("T77").Value = 2.42
$ws.Range("U77").Value = "10/11/2023 11:56"
$ws.Range("V77").Value = "https://www.betexplorer.com/football/thailand/thai-league-1/sukhothai-chonburi/UsyQ5oKG/"

$wb.Save()
